# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect freshly scraped counts (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Row -> (old, new) value for "想去人数" column F, shared by both sheets
# (row numbers differ slightly between sheets because 全部类型 has one
# extra row that 展览 does not).
$updatesByRow = @{
    2  = 185
    4  = 141
    5  = 1294
    6  = 18103
    7  = 360
    8  = 259
    10 = 6820
    11 = 685
    12 = 158
    14 = 110
    17 = 153
    19 = 214
    21 = 654
    22 = 36
    23 = 31
    25 = 272
    26 = 984
    27 = 120
    28 = 5164
    30 = 30
    33 = 12053
    34 = 1279
    35 = 41
    36 = 207
    37 = 277
    38 = 3916
    39 = 300
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updatesByRow.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updatesByRow[$row]
}

# "全部类型" mirrors "展览" for rows 2-28, but rows from 30 onward are
# shifted down by one (row 30 in 展览 corresponds to row 32 here, etc.)
$updatesByRow2 = @{
    2  = 185
    4  = 141
    5  = 1294
    6  = 18103
    7  = 360
    8  = 259
    10 = 6820
    11 = 685
    12 = 158
    14 = 110
    17 = 153
    19 = 214
    21 = 654
    22 = 36
    23 = 31
    25 = 272
    26 = 984
    27 = 120
    28 = 5164
    32 = 30
    35 = 12053
    36 = 1279
    37 = 41
    38 = 207
    39 = 277
    40 = 3916
    41 = 300
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesByRow2.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updatesByRow2[$row]
}
